$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FEB-22")
Write-Output $ws.Name
